# Revert "Increased Slugs and Buckshot damages"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 38 (ammo_12x76_zhekan / Slugs): H38 is a plain value, revert 2.7 -> 2.5
$ws.Range("H38").Value = 2.5

# Row 39 (ammo_12x70_buck / Buckshot): H39 is a formula, revert 9*0.42 -> 9*0.4
$ws.Range("H39").Formula = "=9*0.4"

# Update the active selection to match the saved view state of the sheet
$ws.Activate()
$ws.Range("N22").Select()
